$wb = $excel.ActiveWorkbook

# --- Evcard sheet: close out September with a monthly SUM formula ---
$evcard = $wb.Worksheets.Item("Evcard")
$evcard.Range("C45").Formula = "=SUM(B27:B45)"

# --- Ponycar sheet: append two more trips for end of September ---
$ponycar = $wb.Worksheets.Item("Ponycar")

# Carry over the date/border formatting from the preceding data row
$ponycar.Range("A34:B34").Copy()
$ponycar.Range("A35:B35").PasteSpecial(-4122)
$ponycar.Range("A34:B34").Copy()
$ponycar.Range("A36:B36").PasteSpecial(-4122)

$ponycar.Cells.Item(35, 1).Value = 43373
$ponycar.Cells.Item(35, 2).Value = 4
$ponycar.Cells.Item(36, 1).Value = 43373
$ponycar.Cells.Item(36, 2).Value = 14
$ponycar.Range("C36").Formula = "=SUM(B27:B36)"

# --- Sheet selections: Evcard becomes the active/selected tab ---
$evcard.Activate()
$evcard.Range("F43").Select()
$ponycar.Range("E29").Select()
$wb.Worksheets.Item("SUM").Range("D16").Select()
$evcard.Activate()
